# TC11_CDS_Filter_FileType-TSV.xlsx update
# ------------------------------------------------------------------
# The "ParticipantsTab" row's Cypher query (column B, "query") is
# rewritten: the query that returns Participant/Study/Gender/Samples
# info for files filtered on file_type 'TSV' is updated to use
# OPTIONAL MATCH clauses (picking up participants even when some
# optional relationships are missing) and to sort the collected
# sample ids via apoc.coll.sort(...) before joining them.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Cypher query text that replaces the old "Gender/Samples" query.
$newParticipantsQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['TSV']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id LIMIT 100
"@

# Locate the row whose TabName (column A) is "ParticipantsTab" so the
# edit targets the right row even if row order ever changes, then
# update its "query" value (column B) in place.
$tabCell = $ws.Columns.Item(1).Find("ParticipantsTab")
if ($tabCell -ne $null) {
    $row = $tabCell.Row
} else {
    $row = 2
}

$ws.Cells.Item($row, 2).Value = $newParticipantsQuery

# Keep the window scrolled one row further down (topLeftCell A3 -> A4),
# matching the refreshed view state saved with the workbook.
$win = $excel.ActiveWindow
if ($win -ne $null) {
    try { $win.ScrollRow = 4 } catch {}
    try { $win.TopLeftCell = $ws.Range("A4") } catch {}
}
